$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Nil" to the ReasonToReject column (J) for the two rejected rows (9 and 14)
$ws.Range("J9").Value = "Nil"
$ws.Range("J14").Value = "Nil"

# Move the active cell selection to J14 (matches the sheetView selection change in the diff)
$ws.Range("J14").Select()
